# Fix property_category values on the "building" (建物) and "car" (汽車)
# sheets: both were stamped with the generic "land" category value copied
# from the land sheet; correct them to "building" / "car" respectively.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet - property_category column is I (col 9)
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

# 汽車 (car) sheet - property_category column is H (col 8)
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
$wsCar.Range("H3").Value = "car"
